# Calendar workbook update: record a "P" status in B2 and an (empty)
# value in D2 for the calendar_001 test row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: empty text value (leading apostrophe forces a text entry so the
# cell keeps an explicit, blank string rather than becoming a truly
# empty/no-value cell), written before B2 so the shared-string table
# picks up the blank string ahead of "P".
$ws.Range("D2").Value = "'"
$ws.Range("D2").Style = "Normal"

# B2: status flag "P"
$ws.Range("B2").Value = "P"
$ws.Range("B2").Style = "Normal"
